$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: duplicate of row 2's INN, same ACCOUNT/Name as row 3 (Dima)
$ws.Range("A4").Value = 123456789
$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "00000000000000000000"
$ws.Range("C4").Value = "Dima"
$ws.Range("C4").NumberFormat = "@"

# Row 5: duplicate of row 3's INN, same ACCOUNT/Name as row 3 (Dima)
$ws.Range("A5").Value = 987654321
$ws.Range("A5").NumberFormat = "@"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "00000000000000000000"
$ws.Range("C5").Value = "Dima"
$ws.Range("C5").NumberFormat = "@"

# Row 6: INN only, ACCOUNT/Name left empty (testing "not empty" validation)
$ws.Range("A6").Value = 123456789
$ws.Range("A6").NumberFormat = "@"

$ws.Range("A5").Select()
